$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.726246953010559
$ws.Range("B1").Value = 2.304184913635254
$ws.Range("C1").Value = 2.489066600799561
$ws.Range("D1").Value = 3.209345579147339
$ws.Range("E1").Value = 1.725540637969971
